# "grants por usuario" sheet: the account that rows 59-117 grant privileges
# to changes from "wlima" (row 59, CREATE USER) / "evaldo" (rows 60-117,
# the various GRANTs) to the single new user "yago".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants por usuario")
$ws.Activate() | Out-Null

for ($r = 59; $r -le 117; $r++) {
    $ws.Cells.Item($r, 2).Value = "yago"
}

# Reflect the author's last on-screen scroll position/selection over the
# block of rows that were just updated.
$ws.Range("D59:D117").Select() | Out-Null
